$d = $word.ActiveDocument

# Locate the paragraph that holds the "Ver no Jupiter..." footer line by
# its text content (robust against any index drift).
$target = "Ver no Jupiter Salvar em pdf Salvar em docx"
$jupIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Contains($target)) {
        $jupIndex = $i
    }
}

# The footer block to remove spans: the blank paragraph right before the
# "Ver no Jupiter..." line, that line itself, and the following
# "(c) 2020 ..." credits line. The empty paragraph that used to trail them
# remains, simply closing out what is now the page-break paragraph above.
$blockStart = $d.Paragraphs.Item($jupIndex - 1).Range.Start
$blockEnd = $d.Paragraphs.Item($jupIndex + 1).Range.End

$d.Range($blockStart, $blockEnd).Delete()
